$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number must be forced to Text format
# first, otherwise Excel auto-converts the string into a numeric value (losing
# exact formatting / trailing zeros and introducing floating-point artifacts).
$textCells = @(
    "D5",
    "D8",
    "D10",
    "D11",
    "D12",
    "D19",
    "D21",
    "D22",
    "D23",
    "D26",
    "D28",
    "D29",
    "D31",
    "D32",
    "D33",
    "D35",
    "D40",
    "D41",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D47",
    "D48",
    "D51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.547.25"
$ws.Range("E2").Value = "  -0.02%  "
$ws.Range("D3").Value = "3.133.49"
$ws.Range("E3").Value = "  -1.25%  "
$ws.Range("D5").Value = "572.71"
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  -4.02%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "0.571"
$ws.Range("E8").Value = "  -6.28%  "
$ws.Range("D9").Value = "3.149.84"
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("D10").Value = "0.118"
$ws.Range("E10").Value = "  -2.19%  "
$ws.Range("D11").Value = "6.64"
$ws.Range("E11").Value = "  -2.66%  "
$ws.Range("D12").Value = "0.381"
$ws.Range("E12").Value = "  -2.92%  "
$ws.Range("D13").Value = "3.682.30"
$ws.Range("E13").Value = "  -1.20%  "
$ws.Range("E14").Value = "  -2.25%  "
$ws.Range("D15").Value = "64.590.85"
$ws.Range("E15").Value = "  -0.02%  "
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").Value = "3.148.87"
$ws.Range("E17").Value = "  -1.16%  "
$ws.Range("E18").Value = "  -2.04%  "
$ws.Range("D19").Value = "414.23"
$ws.Range("E19").Value = "  -0.91%  "
$ws.Range("E20").Value = "  -2.00%  "
$ws.Range("D21").Value = "12.45"
$ws.Range("D22").Value = "7.03"
$ws.Range("E22").Value = "  -1.70%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("E25").Value = "  -3.44%  "
$ws.Range("D26").Value = "0.195"
$ws.Range("E26").Value = "  -4.47%  "
$ws.Range("E27").Value = "  -2.78%  "
$ws.Range("D28").Value = "9.02"
$ws.Range("E28").Value = "  +2.50%  "
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("D31").Value = "1.81"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").Value = "21.21"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("D33").Value = "163.27"
$ws.Range("E33").Value = "  +4.50%  "
$ws.Range("E34").Value = "  -4.62%  "
$ws.Range("D35").Value = "6.24"
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("E37").Value = "  -1.18%  "
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").Value = "2.620.74"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").Value = "4.13"
$ws.Range("E40").Value = "  -2.63%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "23.73"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").Value = "38.30"
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("D43").Value = "0.690"
$ws.Range("E43").Value = "  -3.88%  "
$ws.Range("D44").Value = "0.0614"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").Value = "290.51"
$ws.Range("E45").Value = "  -0.93%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "5.26"
$ws.Range("E46").Value = "  -5.66%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "0.0255"
$ws.Range("E47").Value = "  -3.58%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "21.20"
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").Value = "10.49"
$ws.Range("E51").Value = "  +0.66%  "
